$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts the existing columns
# A:F (index/segment name .. PercActivationsFixed) to B:G, carrying their
# values and styles along with them.
$ws.Columns("A:A").Insert()

# Give the new A2:A20 cells the same style as the (now shifted) segment
# name column B2:B20 used to have, i.e. style index "1" (bold, thin box
# border, centered/top aligned) - matching the diff where the numeric
# index cells keep s="1" while the segment-name text cells end up with
# no explicit style.
$ws.Range("B2:B20").Copy()
$ws.Range("A2:A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header for the inserted column - give it the same bold/bordered
# header style ("s=1") as its neighboring header cells before setting
# its text.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B1").Value = "segments"

# Fill column A with the numeric segment index (0-based) for each row,
# and strip the old bold/border style off the segment-name cells in B
# since in the target file those are plain, unstyled cells.
$segments = @("background","back_bumper","back_glass","back_left_door","back_left_light","back_right_door","back_right_light","front_bumper","front_glass","front_left_door","front_left_light","front_right_door","front_right_light","hood","left_mirror","right_mirror","tailgate","trunk","wheel")

for ($i = 0; $i -lt $segments.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).ClearFormats()
    $ws.Cells.Item($row, 2).Value = $segments[$i]
}
